$d = $word.ActiveDocument
$d.AcceptAllRevisions()
